# Se establece la navegacion entre las ventanas login, principal y registro
# Adds a new registration row (row 5) to the "Hoja1" worksheet, mirroring
# the data captured by the new registro window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text fields (new unique values go to sharedStrings automatically)
$ws.Range("A5").Value = "jhasbdjh"
$ws.Range("B5").Value = "asbdjhsa"
$ws.Range("D5").Value = "shadbjsa223"
$ws.Range("F5").Value = "asbdsad@elpdjcn.dsbh.com"
$ws.Range("G5").Value = "No tiene"

# Numeric "Documento" field
$ws.Range("C5").Value = 268361723

# "Contraseña" re-uses the exact same text value as row 2 ("1234"), stored
# as a shared string rather than a number - copy it over so the cell keeps
# its text type instead of being coerced to a number.
$ws.Range("E2").Copy()
$ws.Range("E5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
